$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "arlig innflasjons rate" column (D) entirely: the header (D2)
# and all of the data values (D3:D18) are removed, while leaving the
# (already-applied) number-format styling on D3:D18 intact.
$ws.Range("D2:D18").ClearContents()
